$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.877485834093231
$wsP.Range("D2").Value = 0.9071886033732404
$wsP.Range("E2").Value = 0.08145768007930565
$wsP.Range("F2").Value = 0.6846860865651805

$wsP.Range("B3").Value = 0.877485834093231
$wsP.Range("D3").Value = 0.9669965261226461
$wsP.Range("E3").Value = 0.1798736687265376
$wsP.Range("F3").Value = 0.6115526131609865

$wsP.Range("B4").Value = 0.9071886033732404
$wsP.Range("C4").Value = 0.9669965261226461
$wsP.Range("E4").Value = 0.2635985529173064
$wsP.Range("F4").Value = 0.5441154283809615

$wsP.Range("B5").Value = 0.08145768007930565
$wsP.Range("C5").Value = 0.1798736687265376
$wsP.Range("D5").Value = 0.2635985529173064
$wsP.Range("F5").Value = 0.1325013132211668

$wsP.Range("B6").Value = 0.6846860865651805
$wsP.Range("C6").Value = 0.6115526131609865
$wsP.Range("D6").Value = 0.5441154283809615
$wsP.Range("E6").Value = 0.1325013132211668

# --- Sheet: Estadisticos_DM ---
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = -0.1570009207322911
$wsE.Range("D2").Value = -0.1187138716353724
$wsE.Range("E2").Value = -1.87739438195806
$wsE.Range("F2").Value = 0.4146440492352215

$wsE.Range("B3").Value = 0.1570009207322911
$wsE.Range("D3").Value = -0.04212169814505921
$wsE.Range("E3").Value = -1.411710068743981
$wsE.Range("F3").Value = 0.519463113448287

$wsE.Range("B4").Value = 0.1187138716353724
$wsE.Range("C4").Value = 0.04212169814505921
$wsE.Range("E4").Value = -1.164712882309042
$wsE.Range("F4").Value = 0.6217241858030794

$wsE.Range("B5").Value = 1.87739438195806
$wsE.Range("C5").Value = 1.411710068743981
$wsE.Range("D5").Value = 1.164712882309042
$wsE.Range("F5").Value = 1.597372237729202

$wsE.Range("B6").Value = -0.4146440492352215
$wsE.Range("C6").Value = -0.519463113448287
$wsE.Range("D6").Value = -0.6217241858030794
$wsE.Range("E6").Value = -1.597372237729202
